# Commit message: redefined prefix "ome" instead of ":" (base prefix name).
#
# The workbook's @prefix sheet maps short prefixes to namespace URIs.
# Previously the "base" (default) prefix had an empty name (rendered as
# just ":" when referenced from other sheets, e.g. ":Image", ":pixels").
# This change gives that base prefix an explicit name, "ome", and updates
# every place elsewhere in the workbook that referenced the base prefix
# (":Foo" / "bareCamelCaseName") so it is spelled "ome:Foo" instead.

$wb = $excel.ActiveWorkbook

# --- @prefix sheet: name the base prefix "ome" (was blank) ---------------
$ws = $wb.Worksheets.Item("@prefix")
$ws.Range("A1").Value = "ome"

# --- Image sheet -----------------------------------------------------------
$ws = $wb.Worksheets.Item("Image")
$ws.Range("E3").Value = "ome:pixels"
$ws.Range("B4").Value = "ome:Image"
$ws.Range("E4").Value = "ome:Pixels"

# --- Pixels sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("Pixels")
$ws.Range("D3").Value = "ome:pixelType"
$ws.Range("E3").Value = "ome:dimensionOrder"
$ws.Range("F3").Value = "ome:physicalSizeX"
$ws.Range("G3").Value = "ome:physicalSizeY"
$ws.Range("H3").Value = "ome:sizeC"
$ws.Range("I3").Value = "ome:sizeT"
$ws.Range("J3").Value = "ome:sizeX"
$ws.Range("K3").Value = "ome:sizeY"
$ws.Range("L3").Value = "ome:sizeZ"
$ws.Range("M3").Value = "ome:channel"
$ws.Range("N3").Value = "ome:binData"
$ws.Range("B4").Value = "ome:Pixels"
$ws.Range("D4").Value = "ome:PixelType"
$ws.Range("E4").Value = "ome:DimensionOrder"
$ws.Range("M4").Value = "ome:Channel"
$ws.Range("N4").Value = "ome:BinData"

# --- Channel sheet -----------------------------------------------------------
$ws = $wb.Worksheets.Item("Channel")
$ws.Range("D3").Value = "ome:color"
$ws.Range("B4").Value = "ome:Channel"
$ws.Range("D4").Value = "ome:Color"

# --- Color sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Color")
$ws.Range("B4").Value = "ome:Color"

# --- Binary_Data sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("Binary_Data")
$ws.Range("C3").Value = "ome:bigEndian"
$ws.Range("D3").Value = "ome:data"
$ws.Range("E3").Value = "ome:length"
$ws.Range("B4").Value = "ome:BinData"

# Leave the final selection on the @prefix (first/active) sheet, cell A1,
# matching the workbook's normal "tabSelected" sheet.
$ws = $wb.Worksheets.Item("@prefix")
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
